# Apply the Dec 25, 2023 cryptos-list update (GitHub Actions scraper).
# D/E columns hold Price / Volume(1h) as plain text; values that would
# otherwise be auto-parsed by Excel as a pure number are written with a
# leading apostrophe so they stay text (matching the source data's
# inlineStr/string cell type) instead of silently becoming numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.692.60"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.285.37"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'124.74"
$ws.Range("E5").Value = "  +9.76%  "
$ws.Range("D6").Value = "'266.43"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "'49.04"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'0.0941"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'9.16"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "'15.55"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").Value = "'0.900"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "2.629.47"
$ws.Range("D17").Value = "2.278.13"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "43.737.32"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("E20").Value = "  +3.00%  "
$ws.Range("D21").Value = "'72.46"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "'2.46"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").Value = "'235.85"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'9.64"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "'2.88"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "'11.85"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "'43.13"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").Value = "'3.38"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'173.06"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "'21.75"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'0.0913"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "'5.77"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("E36").Value = "  +5.42%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "'4.15"
$ws.Range("E38").Value = "  +8.76%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").Value = "'2.55"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'76.08"
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'14.15"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'5.87"
$ws.Range("E45").Value = "  -7.51%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "'75.03"
$ws.Range("E47").Value = "  +40.23%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "'8.57"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").Value = "'102.42"
$ws.Range("E51").Value = "  +0.25%  "
